# Apply the "Add data for 2022-06-06" update to the carjacking-by-neighborhood
# workbook. This moves the rolling "through" date from May 28 to May 29 and
# bumps up the counts / adds new cells for the newly-recorded incidents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the worksheet to reflect the new "through" date.
$ws.Name = "Through 2022-05-29"

# 2. Update the header label for the current (rolling) month column (B1).
$ws.Range("B1").Value = "May 2022 (through May 29)"

# 3. Update existing incident counts that changed.
$ws.Range("B2").Value = 12   # Englewood, May 2022
$ws.Range("L2").Value = 5    # Englewood, May 2020
$ws.Range("G3").Value = 10   # Austin, May 2021
$ws.Range("B19").Value = 3   # Little Italy, UIC, May 2022
$ws.Range("B32").Value = 2   # United Center, May 2022

# 4. Fill in newly-populated cells (previously blank).
$ws.Range("AK34").Value = 1  # South Deering, May 2015
$ws.Range("L37").Value = 1   # River North, May 2020
$ws.Range("L38").Value = 1   # Douglas, May 2020
$ws.Range("B51").Value = 1   # Ashburn, May 2022
$ws.Range("AA94").Value = 1  # West Ridge, May 2017
